$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# Sheet "展览" (sheet1)
$ws1.Range("G2").Value = 49.9
$ws1.Range("F4").Value = 13393
$ws1.Range("F9").Value = 126
$ws1.Range("F14").Value = 13370
$ws1.Range("F16").Value = 587
$ws1.Range("F17").Value = 8911
$ws1.Range("F19").Value = 7981
$ws1.Range("F32").Value = 159

# Sheet "演出" (sheet2)
$ws2.Range("F3").Value = 33

# Sheet "全部类型" (sheet4)
$ws4.Range("G3").Value = 49.9
$ws4.Range("F5").Value = 13393
$ws4.Range("F10").Value = 126
$ws4.Range("F15").Value = 13370
$ws4.Range("F17").Value = 587
$ws4.Range("F18").Value = 8911
$ws4.Range("F20").Value = 7981
$ws4.Range("F31").Value = 33
$ws4.Range("F35").Value = 159
